$d = $word.ActiveDocument

# 1. Merge the two runs "There are two circuit design " + "scheme" into a
#    single run with the new, expanded sentence.
$found = $d.Content.Find.Execute(
    "There are two circuit design scheme",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There are two circuit design schemes we and it is hard to choose which one to use.",
    2)
Write-Output "Replace result: $found"

# 2. Insert a new, empty bullet-list paragraph right after that list item
#    (it inherits the ListParagraph/numPr formatting of the paragraph mark
#    it is split off from), pushing the following blank "Normal" paragraph
#    further down.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*There are two circuit design schemes we and it is hard to choose which one to use.*") {
        $end = $p.Range.End
        # Position right before this paragraph's own paragraph mark so the
        # break lands inside the list item, not the following paragraph.
        $splitPoint = $d.Range($end - 1, $end - 1)
        $splitPoint.InsertBefore("`r")
        break
    }
}
